$wb = $excel.ActiveWorkbook

# "Daily Orders" sheet: update order status from NEW to COOKING
$wsOrders = $wb.Worksheets.Item("Daily Orders")
$wsOrders.Range("H2").Value = "COOKING"

# "Summary" sheet: reflect the status change in the summary counts
# (New count goes from 1 to 0, Cooking count goes from 0 to 1)
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0
$wsSummary.Range("C2").Value = 1
